# Apply "some fixs + another bug" changes to the bug-tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-BodyCell($cell, $text) {
    $cell.Value = $text
    $cell.WrapText = $true
    $cell.HorizontalAlignment = -4131
    $cell.VerticalAlignment = -4108
}

# --- Row 8: "Game needs to be more colorful" ---
# Add solution notes in column C and change Programer from "aaron" -> "Col"
Set-BodyCell $ws.Range("C8") "added color to answer button on user selection. `nAdded time delay so users can see it"
Set-BodyCell $ws.Range("D8") "Col"

# --- Row 9: "Answers need to be in random locations" ---
# Add solution notes in column C and change Programer from "aaron" -> "col"
Set-BodyCell $ws.Range("C9") "ramdomised where buttons are "
Set-BodyCell $ws.Range("D9") "col"

# --- Row 11: "Name above menu persisting after logout" ---
# Add solution + programmer info (new bug fix notes)
Set-BodyCell $ws.Range("C11") "added a clicked bool that is checked before a user makes a selection. It is set when they make a selection and reset when the next question is shown."
Set-BodyCell $ws.Range("D11") "col"

# --- Row 12 (new row): multiple answer selections bug ---
Set-BodyCell $ws.Range("A12") "It is possible to make multiple answer selections"
Set-BodyCell $ws.Range("B12") "after a user makes selects an answer they need to blocked from making another selection"
Set-BodyCell $ws.Range("C12") "added a clicked bool that is checked before a user makes a selection. It is set when they make a selection and reset when the next question is shown."
Set-BodyCell $ws.Range("D12") "Col"

# --- Row 13 (new row): duplicate questions bug ---
Set-BodyCell $ws.Range("A13") "duplications of questions are bein shown within the same game"

# Update the active selection / view to match the end state of the edit.
$ws.Activate()
$ws.Range("B13").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
